$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B header / trigger flag for row 1
$ws.Range("B1").Value = "trigger"

# Row 2: stimulus file changes from sad clip to Set2-Vibr, and a paired
# trigger-tagged copy is added in column B
$ws.Range("A2").Value = "Stimuli/Set2-Vibr.wav"
$ws.Range("B2").Value = "Stimuli/trigger_Set2-Vibr.wav"
$ws.Range("B2").WrapText = $True

# Row 3: the old "happy" stimulus row is removed (A3 cleared); B3 keeps
# getting touched/reformatted alongside it, staying empty, picking up the
# same (theme-coloured) font used by the rest of column A/B
$ws.Range("A3").Value = ""
$ws.Range("B3").Value = ""
$ws.Range("B3").WrapText = $True
$ws.Range("B3").Font.ThemeColor = 1
